$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Shift the B:G values of rows 2-10 down into rows 3-11 (row 11's old
# values are discarded), then write the new row 2 values.
for ($r = 10; $r -ge 2; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r + 1, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

$ws.Cells.Item(2, 2).Value2 = -0.02314597604078636
$ws.Cells.Item(2, 3).Value2 = 0.3579920056255013
$ws.Cells.Item(2, 4).Value2 = 0.1782699060034266
$ws.Cells.Item(2, 5).Value2 = 0.4222202103209018
$ws.Cells.Item(2, 6).Value2 = 0.4363822494547141
$ws.Cells.Item(2, 7).Value2 = 15
